$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = 'Última actualización: 18:30:48'
$ws1.Range("A3").Value = 'Total filas: 312'
$arr1 = New-Object 'object[,]' 312,5
$arr1[0,0] = '05:57:04'
$arr1[0,1] = '06:09'
$arr1[0,2] = '10_OLMOS'
$arr1[0,3] = 12
$arr1[0,4] = 'LP1912'
$arr1[1,0] = '05:57:04'
$arr1[1,1] = '06:16'
$arr1[1,2] = '215A_EL PATO'
$arr1[1,3] = 19
$arr1[1,4] = 'LP1912'
$arr1[2,0] = '05:57:04'
$arr1[2,1] = '06:30'
$arr1[2,2] = '23_HERNANDEZ'
$arr1[2,3] = 33
$arr1[2,4] = 'LP1912'
$arr1[3,0] = '05:57:04'
$arr1[3,1] = '06:34'
$arr1[3,2] = '11_ETCHEVERRY'
$arr1[3,3] = 37
$arr1[3,4] = 'LP1912'
$arr1[4,0] = '06:34:35'
$arr1[4,1] = '06:36'
$arr1[4,2] = '11_ETCHEVERRY'
$arr1[4,3] = 2
$arr1[4,4] = 'LP1912'
$arr1[5,0] = '05:57:04'
$arr1[5,1] = '06:39'
$arr1[5,2] = '17X38_ROMERO'
$arr1[5,3] = 42
$arr1[5,4] = 'LP1912'
$arr1[6,0] = '05:57:04'
$arr1[6,1] = '06:41'
$arr1[6,2] = '16_SANTA ANA'
$arr1[6,3] = 44
$arr1[6,4] = 'LP1912'
$arr1[7,0] = '06:16:41'
$arr1[7,1] = '06:56'
$arr1[7,2] = '215A_EL PATO'
$arr1[7,3] = 40
$arr1[7,4] = 'LP1912'
$arr1[8,0] = '05:57:04'
$arr1[8,1] = '06:57'
$arr1[8,2] = '215A_EL PATO'
$arr1[8,3] = 60
$arr1[8,4] = 'LP1912'
$arr1[9,0] = '05:57:04'
$arr1[9,1] = '06:59'
$arr1[9,2] = '225_GOMEZ'
$arr1[9,3] = 62
$arr1[9,4] = 'LP1912'
$arr1[10,0] = '06:16:41'
$arr1[10,1] = '07:15'
$arr1[10,2] = '215C_EL PATO'
$arr1[10,3] = 59
$arr1[10,4] = 'LP1912'
$arr1[11,0] = '05:57:04'
$arr1[11,1] = '07:16'
$arr1[11,2] = '215C_EL PATO'
$arr1[11,3] = 79
$arr1[11,4] = 'LP1912'
$arr1[12,0] = '05:57:04'
$arr1[12,1] = '07:19'
$arr1[12,2] = '14_ABASTO'
$arr1[12,3] = 82
$arr1[12,4] = 'LP1912'
$arr1[13,0] = '06:16:41'
$arr1[13,1] = '07:20'
$arr1[13,2] = '16_SANTA ANA'
$arr1[13,3] = 64
$arr1[13,4] = 'LP1912'
$arr1[14,0] = '05:57:04'
$arr1[14,1] = '07:21'
$arr1[14,2] = '16_SANTA ANA'
$arr1[14,3] = 84
$arr1[14,4] = 'LP1912'
$arr1[15,0] = '06:16:41'
$arr1[15,1] = '07:21'
$arr1[15,2] = '23_HERNANDEZ'
$arr1[15,3] = 65
$arr1[15,4] = 'LP1912'
$arr1[16,0] = '05:57:04'
$arr1[16,1] = '07:22'
$arr1[16,2] = '23_HERNANDEZ'
$arr1[16,3] = 85
$arr1[16,4] = 'LP1912'
$arr1[17,0] = '05:57:04'
$arr1[17,1] = '07:29'
$arr1[17,2] = '17X38_ROMERO'
$arr1[17,3] = 92
$arr1[17,4] = 'LP1912'
$arr1[18,0] = '05:57:04'
$arr1[18,1] = '07:35'
$arr1[18,2] = '10_OLMOS'
$arr1[18,3] = 98
$arr1[18,4] = 'LP1912'
$arr1[19,0] = '06:16:41'
$arr1[19,1] = '07:36'
$arr1[19,2] = '27_EL RETIRO'
$arr1[19,3] = 80
$arr1[19,4] = 'LP1912'
$arr1[20,0] = '05:57:04'
$arr1[20,1] = '07:37'
$arr1[20,2] = '27_EL RETIRO'
$arr1[20,3] = 100
$arr1[20,4] = 'LP1912'
$arr1[21,0] = '06:34:35'
$arr1[21,1] = '07:43'
$arr1[21,2] = '215A_EL PATO'
$arr1[21,3] = 69
$arr1[21,4] = 'LP1912'
$arr1[22,0] = '06:46:20'
$arr1[22,1] = '07:44'
$arr1[22,2] = '215A_EL PATO'
$arr1[22,3] = 58
$arr1[22,4] = 'LP1912'
$arr1[23,0] = '07:49:32'
$arr1[23,1] = '07:49'
$arr1[23,2] = '215A_EL PATO'
$arr1[23,3] = 0
$arr1[23,4] = 'LP1912'
$arr1[24,0] = '05:57:04'
$arr1[24,1] = '07:55'
$arr1[24,2] = '14_ABASTO'
$arr1[24,3] = 118
$arr1[24,4] = 'LP1912'
$arr1[25,0] = '07:56:02'
$arr1[25,1] = '07:59'
$arr1[25,2] = '14_ABASTO'
$arr1[25,3] = 3
$arr1[25,4] = 'LP1912'
$arr1[26,0] = '06:16:41'
$arr1[26,1] = '08:00'
$arr1[26,2] = '17_ROMERO'
$arr1[26,3] = 104
$arr1[26,4] = 'LP1912'
$arr1[27,0] = '06:16:41'
$arr1[27,1] = '08:01'
$arr1[27,2] = '16_SANTA ANA'
$arr1[27,3] = 105
$arr1[27,4] = 'LP1912'
$arr1[28,0] = '06:34:35'
$arr1[28,1] = '08:06'
$arr1[28,2] = '23_HERNANDEZ'
$arr1[28,3] = 92
$arr1[28,4] = 'LP1912'
$arr1[29,0] = '07:56:02'
$arr1[29,1] = '08:07'
$arr1[29,2] = '23_HERNANDEZ'
$arr1[29,3] = 11
$arr1[29,4] = 'LP1912'
$arr1[30,0] = '06:16:41'
$arr1[30,1] = '08:11'
$arr1[30,2] = '10_OLMOS'
$arr1[30,3] = 115
$arr1[30,4] = 'LP1912'
$arr1[31,0] = '06:16:41'
$arr1[31,1] = '08:13'
$arr1[31,2] = '15X38_ABASTO'
$arr1[31,3] = 117
$arr1[31,4] = 'LP1912'
$arr1[32,0] = '06:34:35'
$arr1[32,1] = '08:29'
$arr1[32,2] = '11_ETCHEVERRY'
$arr1[32,3] = 115
$arr1[32,4] = 'LP1912'
$arr1[33,0] = '06:34:35'
$arr1[33,1] = '08:29'
$arr1[33,2] = '15_ABASTO'
$arr1[33,3] = 115
$arr1[33,4] = 'LP1912'
$arr1[34,0] = '08:28:52'
$arr1[34,1] = '08:32'
$arr1[34,2] = '11_ETCHEVERRY'
$arr1[34,3] = 4
$arr1[34,4] = 'LP1912'
$arr1[35,0] = '08:38:24'
$arr1[35,1] = '08:40'
$arr1[35,2] = '10_OLMOS'
$arr1[35,3] = 2
$arr1[35,4] = 'LP1912'
$arr1[36,0] = '06:46:20'
$arr1[36,1] = '08:41'
$arr1[36,2] = '16_P MOR-SANTA ANA'
$arr1[36,3] = 115
$arr1[36,4] = 'LP1912'
$arr1[37,0] = '06:53:44'
$arr1[37,1] = '08:43'
$arr1[37,2] = '215C_EL PATO'
$arr1[37,3] = 110
$arr1[37,4] = 'LP1912'
$arr1[38,0] = '06:46:20'
$arr1[38,1] = '08:44'
$arr1[38,2] = '215C_EL PATO'
$arr1[38,3] = 118
$arr1[38,4] = 'LP1912'
$arr1[39,0] = '07:56:02'
$arr1[39,1] = '08:45'
$arr1[39,2] = '23_HERNANDEZ'
$arr1[39,3] = 49
$arr1[39,4] = 'LP1912'
$arr1[40,0] = '08:45:31'
$arr1[40,1] = '08:45'
$arr1[40,2] = '215C_EL PATO'
$arr1[40,3] = 0
$arr1[40,4] = 'LP1912'
$arr1[41,0] = '07:49:32'
$arr1[41,1] = '08:47'
$arr1[41,2] = '23_HERNANDEZ'
$arr1[41,3] = 58
$arr1[41,4] = 'LP1912'
$arr1[42,0] = '07:38:09'
$arr1[42,1] = '08:51'
$arr1[42,2] = '23_HERNANDEZ'
$arr1[42,3] = 73
$arr1[42,4] = 'LP1912'
$arr1[43,0] = '07:13:03'
$arr1[43,1] = '08:52'
$arr1[43,2] = '23_HERNANDEZ'
$arr1[43,3] = 99
$arr1[43,4] = 'LP1912'
$arr1[44,0] = '08:52:40'
$arr1[44,1] = '08:52'
$arr1[44,2] = '215B_EL PATO'
$arr1[44,3] = 0
$arr1[44,4] = 'LP1912'
$arr1[45,0] = '07:49:32'
$arr1[45,1] = '08:53'
$arr1[45,2] = '215B_EL PATO'
$arr1[45,3] = 64
$arr1[45,4] = 'LP1912'
$arr1[46,0] = '07:13:03'
$arr1[46,1] = '08:54'
$arr1[46,2] = '215B_EL PATO'
$arr1[46,3] = 101
$arr1[46,4] = 'LP1912'
$arr1[47,0] = '08:28:52'
$arr1[47,1] = '08:57'
$arr1[47,2] = '215A_EL PATO'
$arr1[47,3] = 29
$arr1[47,4] = 'LP1912'
$arr1[48,0] = '07:13:03'
$arr1[48,1] = '08:58'
$arr1[48,2] = '215A_EL PATO'
$arr1[48,3] = 105
$arr1[48,4] = 'LP1912'
$arr1[49,0] = '08:28:52'
$arr1[49,1] = '09:04'
$arr1[49,2] = '10_OLMOS'
$arr1[49,3] = 36
$arr1[49,4] = 'LP1912'
$arr1[50,0] = '08:11:18'
$arr1[50,1] = '09:05'
$arr1[50,2] = '10_OLMOS'
$arr1[50,3] = 54
$arr1[50,4] = 'LP1912'
$arr1[51,0] = '07:38:09'
$arr1[51,1] = '09:06'
$arr1[51,2] = '16_SANTA ANA'
$arr1[51,3] = 88
$arr1[51,4] = 'LP1912'
$arr1[52,0] = '07:56:02'
$arr1[52,1] = '09:11'
$arr1[52,2] = '27_EL RETIRO'
$arr1[52,3] = 75
$arr1[52,4] = 'LP1912'
$arr1[53,0] = '07:49:32'
$arr1[53,1] = '09:12'
$arr1[53,2] = '27_EL RETIRO'
$arr1[53,3] = 83
$arr1[53,4] = 'LP1912'
$arr1[54,0] = '07:38:09'
$arr1[54,1] = '09:14'
$arr1[54,2] = '27_EL RETIRO'
$arr1[54,3] = 96
$arr1[54,4] = 'LP1912'
$arr1[55,0] = '08:28:52'
$arr1[55,1] = '09:16'
$arr1[55,2] = '27_EL RETIRO'
$arr1[55,3] = 48
$arr1[55,4] = 'LP1912'
$arr1[56,0] = '08:38:24'
$arr1[56,1] = '09:17'
$arr1[56,2] = '27_EL RETIRO'
$arr1[56,3] = 39
$arr1[56,4] = 'LP1912'
$arr1[57,0] = '07:49:32'
$arr1[57,1] = '09:17'
$arr1[57,2] = '14_ABASTO'
$arr1[57,3] = 88
$arr1[57,4] = 'LP1912'
$arr1[58,0] = '07:38:09'
$arr1[58,1] = '09:18'
$arr1[58,2] = '14_ABASTO'
$arr1[58,3] = 100
$arr1[58,4] = 'LP1912'
$arr1[59,0] = '07:38:09'
$arr1[59,1] = '09:18'
$arr1[59,2] = '15X38_ABASTO'
$arr1[59,3] = 100
$arr1[59,4] = 'LP1912'
$arr1[60,0] = '08:11:18'
$arr1[60,1] = '09:28'
$arr1[60,2] = '23_HERNANDEZ'
$arr1[60,3] = 77
$arr1[60,4] = 'LP1912'
$arr1[61,0] = '08:28:52'
$arr1[61,1] = '09:28'
$arr1[61,2] = '10_OLMOS'
$arr1[61,3] = 60
$arr1[61,4] = 'LP1912'
$arr1[62,0] = '07:38:09'
$arr1[62,1] = '09:29'
$arr1[62,2] = '10_OLMOS'
$arr1[62,3] = 111
$arr1[62,4] = 'LP1912'
$arr1[63,0] = '08:11:18'
$arr1[63,1] = '09:31'
$arr1[63,2] = '16_SANTA ANA'
$arr1[63,3] = 80
$arr1[63,4] = 'LP1912'
$arr1[64,0] = '08:28:52'
$arr1[64,1] = '09:32'
$arr1[64,2] = '23_HERNANDEZ'
$arr1[64,3] = 64
$arr1[64,4] = 'LP1912'
$arr1[65,0] = '08:38:24'
$arr1[65,1] = '09:34'
$arr1[65,2] = '23_HERNANDEZ'
$arr1[65,3] = 56
$arr1[65,4] = 'LP1912'
$arr1[66,0] = '08:52:40'
$arr1[66,1] = '09:36'
$arr1[66,2] = '23_HERNANDEZ'
$arr1[66,3] = 44
$arr1[66,4] = 'LP1912'
$arr1[67,0] = '07:49:32'
$arr1[67,1] = '09:39'
$arr1[67,2] = '15_ABASTO'
$arr1[67,3] = 110
$arr1[67,4] = 'LP1912'
$arr1[68,0] = '07:49:32'
$arr1[68,1] = '09:41'
$arr1[68,2] = '11_ETCHEVERRY'
$arr1[68,3] = 112
$arr1[68,4] = 'LP1912'
$arr1[69,0] = '07:56:02'
$arr1[69,1] = '09:42'
$arr1[69,2] = '11_ETCHEVERRY'
$arr1[69,3] = 106
$arr1[69,4] = 'LP1912'
$arr1[70,0] = '07:49:32'
$arr1[70,1] = '09:43'
$arr1[70,2] = '16_P MOR-SANTA ANA'
$arr1[70,3] = 114
$arr1[70,4] = 'LP1912'
$arr1[71,0] = '08:11:18'
$arr1[71,1] = '09:53'
$arr1[71,2] = '10_OLMOS'
$arr1[71,3] = 102
$arr1[71,4] = 'LP1912'
$arr1[72,0] = '08:28:52'
$arr1[72,1] = '09:58'
$arr1[72,2] = '215C_EL PATO'
$arr1[72,3] = 90
$arr1[72,4] = 'LP1912'
$arr1[73,0] = '08:11:18'
$arr1[73,1] = '09:59'
$arr1[73,2] = '215C_EL PATO'
$arr1[73,3] = 108
$arr1[73,4] = 'LP1912'
$arr1[74,0] = '08:28:52'
$arr1[74,1] = '10:05'
$arr1[74,2] = '14_ABASTO'
$arr1[74,3] = 97
$arr1[74,4] = 'LP1912'
$arr1[75,0] = '08:11:18'
$arr1[75,1] = '10:06'
$arr1[75,2] = '14_ABASTO'
$arr1[75,3] = 115
$arr1[75,4] = 'LP1912'
$arr1[76,0] = '08:28:52'
$arr1[76,1] = '10:13'
$arr1[76,2] = '17X38_ROMERO'
$arr1[76,3] = 105
$arr1[76,4] = 'LP1912'
$arr1[77,0] = '09:22:34'
$arr1[77,1] = '10:21'
$arr1[77,2] = '23_HERNANDEZ'
$arr1[77,3] = 59
$arr1[77,4] = 'LP1912'
$arr1[78,0] = '10:04:30'
$arr1[78,1] = '10:22'
$arr1[78,2] = '23_HERNANDEZ'
$arr1[78,3] = 18
$arr1[78,4] = 'LP1912'
$arr1[79,0] = '09:22:34'
$arr1[79,1] = '10:25'
$arr1[79,2] = '16_SANTA ANA'
$arr1[79,3] = 63
$arr1[79,4] = 'LP1912'
$arr1[80,0] = '10:04:30'
$arr1[80,1] = '10:29'
$arr1[80,2] = '14_ABASTO'
$arr1[80,3] = 25
$arr1[80,4] = 'LP1912'
$arr1[81,0] = '08:38:24'
$arr1[81,1] = '10:29'
$arr1[81,2] = '15_ABASTO'
$arr1[81,3] = 111
$arr1[81,4] = 'LP1912'
$arr1[82,0] = '10:36:50'
$arr1[82,1] = '10:43'
$arr1[82,2] = '11X44_ETCHEVERRY'
$arr1[82,3] = 7
$arr1[82,4] = 'LP1912'
$arr1[83,0] = '08:45:31'
$arr1[83,1] = '10:44'
$arr1[83,2] = '11X44_ETCHEVERRY'
$arr1[83,3] = 119
$arr1[83,4] = 'LP1912'
$arr1[84,0] = '08:52:40'
$arr1[84,1] = '10:46'
$arr1[84,2] = '15_P INDUSTRIAL'
$arr1[84,3] = 114
$arr1[84,4] = 'LP1912'
$arr1[85,0] = '09:22:34'
$arr1[85,1] = '10:53'
$arr1[85,2] = '27_EL RETIRO'
$arr1[85,3] = 91
$arr1[85,4] = 'LP1912'
$arr1[86,0] = '10:36:50'
$arr1[86,1] = '10:55'
$arr1[86,2] = '16_SANTA ANA'
$arr1[86,3] = 19
$arr1[86,4] = 'LP1912'
$arr1[87,0] = '10:04:30'
$arr1[87,1] = '10:56'
$arr1[87,2] = '27_EL RETIRO'
$arr1[87,3] = 52
$arr1[87,4] = 'LP1912'
$arr1[88,0] = '10:56:15'
$arr1[88,1] = '10:57'
$arr1[88,2] = '27_EL RETIRO'
$arr1[88,3] = 1
$arr1[88,4] = 'LP1912'
$arr1[89,0] = '09:22:34'
$arr1[89,1] = '10:57'
$arr1[89,2] = '10_OLMOS'
$arr1[89,3] = 95
$arr1[89,4] = 'LP1912'
$arr1[90,0] = '10:04:30'
$arr1[90,1] = '10:59'
$arr1[90,2] = '10_OLMOS'
$arr1[90,3] = 55
$arr1[90,4] = 'LP1912'
$arr1[91,0] = '09:22:34'
$arr1[91,1] = '11:01'
$arr1[91,2] = '81_EL PELIGRO'
$arr1[91,3] = 99
$arr1[91,4] = 'LP1912'
$arr1[92,0] = '10:04:30'
$arr1[92,1] = '11:03'
$arr1[92,2] = '23_HERNANDEZ'
$arr1[92,3] = 59
$arr1[92,4] = 'LP1912'
$arr1[93,0] = '10:36:50'
$arr1[93,1] = '11:06'
$arr1[93,2] = '23_HERNANDEZ'
$arr1[93,3] = 30
$arr1[93,4] = 'LP1912'
$arr1[94,0] = '09:22:34'
$arr1[94,1] = '11:10'
$arr1[94,2] = '16_P MOR-SANTA ANA'
$arr1[94,3] = 108
$arr1[94,4] = 'LP1912'
$arr1[95,0] = '09:22:34'
$arr1[95,1] = '11:14'
$arr1[95,2] = '14_ABASTO'
$arr1[95,3] = 112
$arr1[95,4] = 'LP1912'
$arr1[96,0] = '10:56:15'
$arr1[96,1] = '11:15'
$arr1[96,2] = '14_ABASTO'
$arr1[96,3] = 19
$arr1[96,4] = 'LP1912'
$arr1[97,0] = '09:22:34'
$arr1[97,1] = '11:15'
$arr1[97,2] = '15X38_ABASTO'
$arr1[97,3] = 113
$arr1[97,4] = 'LP1912'
$arr1[98,0] = '11:13:15'
$arr1[98,1] = '11:17'
$arr1[98,2] = '14_ABASTO'
$arr1[98,3] = 4
$arr1[98,4] = 'LP1912'
$arr1[99,0] = '10:36:50'
$arr1[99,1] = '11:25'
$arr1[99,2] = '16_SANTA ANA'
$arr1[99,3] = 49
$arr1[99,4] = 'LP1912'
$arr1[100,0] = '10:04:30'
$arr1[100,1] = '11:29'
$arr1[100,2] = '10_OLMOS'
$arr1[100,3] = 85
$arr1[100,4] = 'LP1912'
$arr1[101,0] = '10:04:30'
$arr1[101,1] = '11:29'
$arr1[101,2] = '16_SANTA ANA'
$arr1[101,3] = 85
$arr1[101,4] = 'LP1912'
$arr1[102,0] = '10:36:50'
$arr1[102,1] = '11:30'
$arr1[102,2] = '215C_EL PATO'
$arr1[102,3] = 54
$arr1[102,4] = 'LP1912'
$arr1[103,0] = '10:04:30'
$arr1[103,1] = '11:31'
$arr1[103,2] = '215C_EL PATO'
$arr1[103,3] = 87
$arr1[103,4] = 'LP1912'
$arr1[104,0] = '10:04:30'
$arr1[104,1] = '11:41'
$arr1[104,2] = '215B_EL PATO'
$arr1[104,3] = 97
$arr1[104,4] = 'LP1912'
$arr1[105,0] = '10:56:15'
$arr1[105,1] = '11:42'
$arr1[105,2] = '215B_EL PATO'
$arr1[105,3] = 46
$arr1[105,4] = 'LP1912'
$arr1[106,0] = '10:04:30'
$arr1[106,1] = '11:45'
$arr1[106,2] = '15X38_ABASTO'
$arr1[106,3] = 101
$arr1[106,4] = 'LP1912'
$arr1[107,0] = '11:46:32'
$arr1[107,1] = '11:46'
$arr1[107,2] = '15X38_ABASTO'
$arr1[107,3] = 0
$arr1[107,4] = 'LP1912'
$arr1[108,0] = '10:56:15'
$arr1[108,1] = '11:46'
$arr1[108,2] = '23_HERNANDEZ'
$arr1[108,3] = 50
$arr1[108,4] = 'LP1912'
$arr1[109,0] = '10:49:38'
$arr1[109,1] = '11:47'
$arr1[109,2] = '23_HERNANDEZ'
$arr1[109,3] = 58
$arr1[109,4] = 'LP1912'
$arr1[110,0] = '10:36:50'
$arr1[110,1] = '11:48'
$arr1[110,2] = '23_HERNANDEZ'
$arr1[110,3] = 72
$arr1[110,4] = 'LP1912'
$arr1[111,0] = '11:13:15'
$arr1[111,1] = '11:51'
$arr1[111,2] = '23_HERNANDEZ'
$arr1[111,3] = 38
$arr1[111,4] = 'LP1912'
$arr1[112,0] = '11:33:52'
$arr1[112,1] = '11:52'
$arr1[112,2] = '23_HERNANDEZ'
$arr1[112,3] = 19
$arr1[112,4] = 'LP1912'
$arr1[113,0] = '10:36:50'
$arr1[113,1] = '11:52'
$arr1[113,2] = '225_GOMEZ'
$arr1[113,3] = 76
$arr1[113,4] = 'LP1912'
$arr1[114,0] = '10:04:30'
$arr1[114,1] = '11:53'
$arr1[114,2] = '225_GOMEZ'
$arr1[114,3] = 109
$arr1[114,4] = 'LP1912'
$arr1[115,0] = '10:04:30'
$arr1[115,1] = '11:58'
$arr1[115,2] = '17_ROMERO'
$arr1[115,3] = 114
$arr1[115,4] = 'LP1912'
$arr1[116,0] = '10:36:50'
$arr1[116,1] = '12:05'
$arr1[116,2] = '11_ETCHEVERRY'
$arr1[116,3] = 89
$arr1[116,4] = 'LP1912'
$arr1[117,0] = '10:56:15'
$arr1[117,1] = '12:06'
$arr1[117,2] = '11_ETCHEVERRY'
$arr1[117,3] = 70
$arr1[117,4] = 'LP1912'
$arr1[118,0] = '10:36:50'
$arr1[118,1] = '12:10'
$arr1[118,2] = '16_P MOR-SANTA ANA'
$arr1[118,3] = 94
$arr1[118,4] = 'LP1912'
$arr1[119,0] = '10:36:50'
$arr1[119,1] = '12:10'
$arr1[119,2] = '15_ABASTO'
$arr1[119,3] = 94
$arr1[119,4] = 'LP1912'
$arr1[120,0] = '12:11:21'
$arr1[120,1] = '12:11'
$arr1[120,2] = '16_P MOR-SANTA ANA'
$arr1[120,3] = 0
$arr1[120,4] = 'LP1912'
$arr1[121,0] = '12:11:21'
$arr1[121,1] = '12:12'
$arr1[121,2] = '15_ABASTO'
$arr1[121,3] = 1
$arr1[121,4] = 'LP1912'
$arr1[122,0] = '11:33:52'
$arr1[122,1] = '12:16'
$arr1[122,2] = '10_OLMOS'
$arr1[122,3] = 43
$arr1[122,4] = 'LP1912'
$arr1[123,0] = '11:13:15'
$arr1[123,1] = '12:17'
$arr1[123,2] = '10_OLMOS'
$arr1[123,3] = 64
$arr1[123,4] = 'LP1912'
$arr1[124,0] = '10:36:50'
$arr1[124,1] = '12:21'
$arr1[124,2] = '215C_EL PATO'
$arr1[124,3] = 105
$arr1[124,4] = 'LP1912'
$arr1[125,0] = '10:56:15'
$arr1[125,1] = '12:22'
$arr1[125,2] = '215C_EL PATO'
$arr1[125,3] = 86
$arr1[125,4] = 'LP1912'
$arr1[126,0] = '11:13:15'
$arr1[126,1] = '12:29'
$arr1[126,2] = '23_HERNANDEZ'
$arr1[126,3] = 76
$arr1[126,4] = 'LP1912'
$arr1[127,0] = '11:33:52'
$arr1[127,1] = '12:32'
$arr1[127,2] = '23_HERNANDEZ'
$arr1[127,3] = 59
$arr1[127,4] = 'LP1912'
$arr1[128,0] = '10:36:50'
$arr1[128,1] = '12:32'
$arr1[128,2] = '14_ABASTO'
$arr1[128,3] = 116
$arr1[128,4] = 'LP1912'
$arr1[129,0] = '10:56:15'
$arr1[129,1] = '12:33'
$arr1[129,2] = '27_EL RETIRO'
$arr1[129,3] = 97
$arr1[129,4] = 'LP1912'
$arr1[130,0] = '10:56:15'
$arr1[130,1] = '12:33'
$arr1[130,2] = '14_ABASTO'
$arr1[130,3] = 97
$arr1[130,4] = 'LP1912'
$arr1[131,0] = '10:36:50'
$arr1[131,1] = '12:34'
$arr1[131,2] = '15_ABASTO'
$arr1[131,3] = 118
$arr1[131,4] = 'LP1912'
$arr1[132,0] = '11:46:32'
$arr1[132,1] = '12:34'
$arr1[132,2] = '23_HERNANDEZ'
$arr1[132,3] = 48
$arr1[132,4] = 'LP1912'
$arr1[133,0] = '11:53:44'
$arr1[133,1] = '12:36'
$arr1[133,2] = '23_HERNANDEZ'
$arr1[133,3] = 43
$arr1[133,4] = 'LP1912'
$arr1[134,0] = '10:49:38'
$arr1[134,1] = '12:36'
$arr1[134,2] = '27_EL RETIRO'
$arr1[134,3] = 107
$arr1[134,4] = 'LP1912'
$arr1[135,0] = '12:11:21'
$arr1[135,1] = '12:37'
$arr1[135,2] = '27_EL RETIRO'
$arr1[135,3] = 26
$arr1[135,4] = 'LP1912'
$arr1[136,0] = '12:11:21'
$arr1[136,1] = '12:37'
$arr1[136,2] = '23_HERNANDEZ'
$arr1[136,3] = 26
$arr1[136,4] = 'LP1912'
$arr1[137,0] = '11:33:52'
$arr1[137,1] = '12:47'
$arr1[137,2] = '14_ABASTO'
$arr1[137,3] = 74
$arr1[137,4] = 'LP1912'
$arr1[138,0] = '10:49:38'
$arr1[138,1] = '12:48'
$arr1[138,2] = '16_SANTA ANA'
$arr1[138,3] = 119
$arr1[138,4] = 'LP1912'
$arr1[139,0] = '11:33:52'
$arr1[139,1] = '12:48'
$arr1[139,2] = '15X38_ABASTO'
$arr1[139,3] = 75
$arr1[139,4] = 'LP1912'
$arr1[140,0] = '11:33:52'
$arr1[140,1] = '13:02'
$arr1[140,2] = '11_ETCHEVERRY'
$arr1[140,3] = 89
$arr1[140,4] = 'LP1912'
$arr1[141,0] = '11:33:52'
$arr1[141,1] = '13:03'
$arr1[141,2] = '215C_EL PATO'
$arr1[141,3] = 90
$arr1[141,4] = 'LP1912'
$arr1[142,0] = '11:13:15'
$arr1[142,1] = '13:03'
$arr1[142,2] = '11_ETCHEVERRY'
$arr1[142,3] = 110
$arr1[142,4] = 'LP1912'
$arr1[143,0] = '11:46:32'
$arr1[143,1] = '13:04'
$arr1[143,2] = '215C_EL PATO'
$arr1[143,3] = 78
$arr1[143,4] = 'LP1912'
$arr1[144,0] = '11:33:52'
$arr1[144,1] = '13:13'
$arr1[144,2] = '16_SANTA ANA'
$arr1[144,3] = 100
$arr1[144,4] = 'LP1912'
$arr1[145,0] = '11:33:52'
$arr1[145,1] = '13:17'
$arr1[145,2] = '10_OLMOS'
$arr1[145,3] = 104
$arr1[145,4] = 'LP1912'
$arr1[146,0] = '12:46:07'
$arr1[146,1] = '13:19'
$arr1[146,2] = '15_ABASTO'
$arr1[146,3] = 33
$arr1[146,4] = 'LP1912'
$arr1[147,0] = '11:53:44'
$arr1[147,1] = '13:21'
$arr1[147,2] = '23_HERNANDEZ'
$arr1[147,3] = 88
$arr1[147,4] = 'LP1912'
$arr1[148,0] = '12:46:07'
$arr1[148,1] = '13:22'
$arr1[148,2] = '23_HERNANDEZ'
$arr1[148,3] = 36
$arr1[148,4] = 'LP1912'
$arr1[149,0] = '12:33:02'
$arr1[149,1] = '13:23'
$arr1[149,2] = '23_HERNANDEZ'
$arr1[149,3] = 50
$arr1[149,4] = 'LP1912'
$arr1[150,0] = '12:11:21'
$arr1[150,1] = '13:24'
$arr1[150,2] = '23_HERNANDEZ'
$arr1[150,3] = 73
$arr1[150,4] = 'LP1912'
$arr1[151,0] = '11:33:52'
$arr1[151,1] = '13:25'
$arr1[151,2] = '16_P MOR-SANTA ANA'
$arr1[151,3] = 112
$arr1[151,4] = 'LP1912'
$arr1[152,0] = '11:53:44'
$arr1[152,1] = '13:32'
$arr1[152,2] = '215A_EL PATO'
$arr1[152,3] = 99
$arr1[152,4] = 'LP1912'
$arr1[153,0] = '12:11:21'
$arr1[153,1] = '13:32'
$arr1[153,2] = '14_ABASTO'
$arr1[153,3] = 81
$arr1[153,4] = 'LP1912'
$arr1[154,0] = '12:33:02'
$arr1[154,1] = '13:33'
$arr1[154,2] = '14_ABASTO'
$arr1[154,3] = 60
$arr1[154,4] = 'LP1912'
$arr1[155,0] = '11:46:32'
$arr1[155,1] = '13:33'
$arr1[155,2] = '215A_EL PATO'
$arr1[155,3] = 107
$arr1[155,4] = 'LP1912'
$arr1[156,0] = '13:41:21'
$arr1[156,1] = '13:42'
$arr1[156,2] = '81_EL PELIGRO'
$arr1[156,3] = 1
$arr1[156,4] = 'LP1912'
$arr1[157,0] = '13:41:21'
$arr1[157,1] = '13:44'
$arr1[157,2] = '225_GOMEZ'
$arr1[157,3] = 3
$arr1[157,4] = 'LP1912'
$arr1[158,0] = '11:53:44'
$arr1[158,1] = '13:47'
$arr1[158,2] = '225_GOMEZ'
$arr1[158,3] = 114
$arr1[158,4] = 'LP1912'
$arr1[159,0] = '12:33:02'
$arr1[159,1] = '13:54'
$arr1[159,2] = '15_ABASTO'
$arr1[159,3] = 81
$arr1[159,4] = 'LP1912'
$arr1[160,0] = '13:55:43'
$arr1[160,1] = '13:55'
$arr1[160,2] = '15_ABASTO'
$arr1[160,3] = 0
$arr1[160,4] = 'LP1912'
$arr1[161,0] = '13:55:43'
$arr1[161,1] = '13:56'
$arr1[161,2] = '81_EL PELIGRO'
$arr1[161,3] = 1
$arr1[161,4] = 'LP1912'
$arr1[162,0] = '13:55:43'
$arr1[162,1] = '13:58'
$arr1[162,2] = '10_OLMOS'
$arr1[162,3] = 3
$arr1[162,4] = 'LP1912'
$arr1[163,0] = '13:14:31'
$arr1[163,1] = '14:02'
$arr1[163,2] = '16_SANTA ANA'
$arr1[163,3] = 48
$arr1[163,4] = 'LP1912'
$arr1[164,0] = '12:46:07'
$arr1[164,1] = '14:02'
$arr1[164,2] = '23_HERNANDEZ'
$arr1[164,3] = 76
$arr1[164,4] = 'LP1912'
$arr1[165,0] = '12:33:02'
$arr1[165,1] = '14:02'
$arr1[165,2] = '10_OLMOS'
$arr1[165,3] = 89
$arr1[165,4] = 'LP1912'
$arr1[166,0] = '13:14:31'
$arr1[166,1] = '14:05'
$arr1[166,2] = '23_HERNANDEZ'
$arr1[166,3] = 51
$arr1[166,4] = 'LP1912'
$arr1[167,0] = '13:41:21'
$arr1[167,1] = '14:06'
$arr1[167,2] = '23_HERNANDEZ'
$arr1[167,3] = 25
$arr1[167,4] = 'LP1912'
$arr1[168,0] = '12:46:07'
$arr1[168,1] = '14:08'
$arr1[168,2] = '16_SANTA ANA'
$arr1[168,3] = 82
$arr1[168,4] = 'LP1912'
$arr1[169,0] = '12:53:26'
$arr1[169,1] = '14:09'
$arr1[169,2] = '23_HERNANDEZ'
$arr1[169,3] = 76
$arr1[169,4] = 'LP1912'
$arr1[170,0] = '13:41:21'
$arr1[170,1] = '14:14'
$arr1[170,2] = '15_ABASTO'
$arr1[170,3] = 33
$arr1[170,4] = 'LP1912'
$arr1[171,0] = '12:53:26'
$arr1[171,1] = '14:16'
$arr1[171,2] = '27_EL RETIRO'
$arr1[171,3] = 83
$arr1[171,4] = 'LP1912'
$arr1[172,0] = '12:53:26'
$arr1[172,1] = '14:17'
$arr1[172,2] = '11_ETCHEVERRY'
$arr1[172,3] = 84
$arr1[172,4] = 'LP1912'
$arr1[173,0] = '12:33:02'
$arr1[173,1] = '14:17'
$arr1[173,2] = '27_EL RETIRO'
$arr1[173,3] = 104
$arr1[173,4] = 'LP1912'
$arr1[174,0] = '12:33:02'
$arr1[174,1] = '14:18'
$arr1[174,2] = '11_ETCHEVERRY'
$arr1[174,3] = 105
$arr1[174,4] = 'LP1912'
$arr1[175,0] = '12:53:26'
$arr1[175,1] = '14:27'
$arr1[175,2] = '16_SANTA ANA'
$arr1[175,3] = 94
$arr1[175,4] = 'LP1912'
$arr1[176,0] = '14:11:28'
$arr1[176,1] = '14:29'
$arr1[176,2] = '10_OLMOS'
$arr1[176,3] = 18
$arr1[176,4] = 'LP1912'
$arr1[177,0] = '12:33:02'
$arr1[177,1] = '14:32'
$arr1[177,2] = '14X44_ABASTO'
$arr1[177,3] = 119
$arr1[177,4] = 'LP1912'
$arr1[178,0] = '13:55:43'
$arr1[178,1] = '14:33'
$arr1[178,2] = '215C_EL PATO'
$arr1[178,3] = 38
$arr1[178,4] = 'LP1912'
$arr1[179,0] = '14:32:44'
$arr1[179,1] = '14:33'
$arr1[179,2] = '14X44_ABASTO'
$arr1[179,3] = 1
$arr1[179,4] = 'LP1912'
$arr1[180,0] = '12:46:07'
$arr1[180,1] = '14:34'
$arr1[180,2] = '215C_EL PATO'
$arr1[180,3] = 108
$arr1[180,4] = 'LP1912'
$arr1[181,0] = '12:46:07'
$arr1[181,1] = '14:39'
$arr1[181,2] = '16_P MOR-SANTA ANA'
$arr1[181,3] = 113
$arr1[181,4] = 'LP1912'
$arr1[182,0] = '12:53:26'
$arr1[182,1] = '14:47'
$arr1[182,2] = '215B_EL PATO'
$arr1[182,3] = 114
$arr1[182,4] = 'LP1912'
$arr1[183,0] = '14:46:12'
$arr1[183,1] = '14:48'
$arr1[183,2] = '215B_EL PATO'
$arr1[183,3] = 2
$arr1[183,4] = 'LP1912'
$arr1[184,0] = '13:41:21'
$arr1[184,1] = '14:51'
$arr1[184,2] = '23_HERNANDEZ'
$arr1[184,3] = 70
$arr1[184,4] = 'LP1912'
$arr1[185,0] = '13:55:43'
$arr1[185,1] = '14:51'
$arr1[185,2] = '16_SANTA ANA'
$arr1[185,3] = 56
$arr1[185,4] = 'LP1912'
$arr1[186,0] = '13:55:43'
$arr1[186,1] = '14:53'
$arr1[186,2] = '215_EL PELIGRO'
$arr1[186,3] = 58
$arr1[186,4] = 'LP1912'
$arr1[187,0] = '13:14:31'
$arr1[187,1] = '14:54'
$arr1[187,2] = '215_EL PELIGRO'
$arr1[187,3] = 100
$arr1[187,4] = 'LP1912'
$arr1[188,0] = '13:14:31'
$arr1[188,1] = '15:02'
$arr1[188,2] = '10_OLMOS'
$arr1[188,3] = 108
$arr1[188,4] = 'LP1912'
$arr1[189,0] = '13:14:31'
$arr1[189,1] = '15:13'
$arr1[189,2] = '17X38_ROMERO'
$arr1[189,3] = 119
$arr1[189,4] = 'LP1912'
$arr1[190,0] = '14:32:44'
$arr1[190,1] = '15:16'
$arr1[190,2] = '16_SANTA ANA'
$arr1[190,3] = 44
$arr1[190,4] = 'LP1912'
$arr1[191,0] = '13:55:43'
$arr1[191,1] = '15:17'
$arr1[191,2] = '14_ABASTO'
$arr1[191,3] = 82
$arr1[191,4] = 'LP1912'
$arr1[192,0] = '13:41:21'
$arr1[192,1] = '15:18'
$arr1[192,2] = '14_ABASTO'
$arr1[192,3] = 97
$arr1[192,4] = 'LP1912'
$arr1[193,0] = '13:55:43'
$arr1[193,1] = '15:33'
$arr1[193,2] = '215C_EL PATO'
$arr1[193,3] = 98
$arr1[193,4] = 'LP1912'
$arr1[194,0] = '13:41:21'
$arr1[194,1] = '15:34'
$arr1[194,2] = '215C_EL PATO'
$arr1[194,3] = 113
$arr1[194,4] = 'LP1912'
$arr1[195,0] = '14:11:28'
$arr1[195,1] = '15:36'
$arr1[195,2] = '23_HERNANDEZ'
$arr1[195,3] = 85
$arr1[195,4] = 'LP1912'
$arr1[196,0] = '13:55:43'
$arr1[196,1] = '15:41'
$arr1[196,2] = '11_ETCHEVERRY'
$arr1[196,3] = 106
$arr1[196,4] = 'LP1912'
$arr1[197,0] = '14:46:12'
$arr1[197,1] = '15:42'
$arr1[197,2] = '11_ETCHEVERRY'
$arr1[197,3] = 56
$arr1[197,4] = 'LP1912'
$arr1[198,0] = '14:32:44'
$arr1[198,1] = '15:53'
$arr1[198,2] = '10_OLMOS'
$arr1[198,3] = 81
$arr1[198,4] = 'LP1912'
$arr1[199,0] = '13:55:43'
$arr1[199,1] = '15:53'
$arr1[199,2] = '16_P MOR-SANTA ANA'
$arr1[199,3] = 118
$arr1[199,4] = 'LP1912'
$arr1[200,0] = '13:55:43'
$arr1[200,1] = '15:53'
$arr1[200,2] = '15X38_ABASTO'
$arr1[200,3] = 118
$arr1[200,4] = 'LP1912'
$arr1[201,0] = '14:46:12'
$arr1[201,1] = '15:54'
$arr1[201,2] = '27_EL RETIRO'
$arr1[201,3] = 68
$arr1[201,4] = 'LP1912'
$arr1[202,0] = '14:11:28'
$arr1[202,1] = '15:56'
$arr1[202,2] = '27_EL RETIRO'
$arr1[202,3] = 105
$arr1[202,4] = 'LP1912'
$arr1[203,0] = '14:53:29'
$arr1[203,1] = '16:02'
$arr1[203,2] = '16_SANTA ANA'
$arr1[203,3] = 69
$arr1[203,4] = 'LP1912'
$arr1[204,0] = '14:11:28'
$arr1[204,1] = '16:05'
$arr1[204,2] = '14_ABASTO'
$arr1[204,3] = 114
$arr1[204,4] = 'LP1912'
$arr1[205,0] = '15:16:46'
$arr1[205,1] = '16:05'
$arr1[205,2] = '16_SANTA ANA'
$arr1[205,3] = 49
$arr1[205,4] = 'LP1912'
$arr1[206,0] = '14:46:12'
$arr1[206,1] = '16:06'
$arr1[206,2] = '14_ABASTO'
$arr1[206,3] = 80
$arr1[206,4] = 'LP1912'
$arr1[207,0] = '15:56:56'
$arr1[207,1] = '16:13'
$arr1[207,2] = '17_ROMERO'
$arr1[207,3] = 17
$arr1[207,4] = 'LP1912'
$arr1[208,0] = '14:32:44'
$arr1[208,1] = '16:14'
$arr1[208,2] = '17_ROMERO'
$arr1[208,3] = 102
$arr1[208,4] = 'LP1912'
$arr1[209,0] = '15:56:56'
$arr1[209,1] = '16:16'
$arr1[209,2] = '10_OLMOS'
$arr1[209,3] = 20
$arr1[209,4] = 'LP1912'
$arr1[210,0] = '14:46:12'
$arr1[210,1] = '16:17'
$arr1[210,2] = '10_OLMOS'
$arr1[210,3] = 91
$arr1[210,4] = 'LP1912'
$arr1[211,0] = '14:32:44'
$arr1[211,1] = '16:21'
$arr1[211,2] = '23_HERNANDEZ'
$arr1[211,3] = 109
$arr1[211,4] = 'LP1912'
$arr1[212,0] = '15:16:46'
$arr1[212,1] = '16:22'
$arr1[212,2] = '23_HERNANDEZ'
$arr1[212,3] = 66
$arr1[212,4] = 'LP1912'
$arr1[213,0] = '15:44:42'
$arr1[213,1] = '16:29'
$arr1[213,2] = '14_ABASTO'
$arr1[213,3] = 45
$arr1[213,4] = 'LP1912'
$arr1[214,0] = '15:16:46'
$arr1[214,1] = '16:30'
$arr1[214,2] = '16_SANTA ANA'
$arr1[214,3] = 74
$arr1[214,4] = 'LP1912'
$arr1[215,0] = '16:12:06'
$arr1[215,1] = '16:30'
$arr1[215,2] = '14_ABASTO'
$arr1[215,3] = 18
$arr1[215,4] = 'LP1912'
$arr1[216,0] = '16:28:21'
$arr1[216,1] = '16:32'
$arr1[216,2] = '14_ABASTO'
$arr1[216,3] = 4
$arr1[216,4] = 'LP1912'
$arr1[217,0] = '15:56:56'
$arr1[217,1] = '16:33'
$arr1[217,2] = '83_ALUAR'
$arr1[217,3] = 37
$arr1[217,4] = 'LP1912'
$arr1[218,0] = '14:46:12'
$arr1[218,1] = '16:34'
$arr1[218,2] = '83_ALUAR'
$arr1[218,3] = 108
$arr1[218,4] = 'LP1912'
$arr1[219,0] = '15:56:56'
$arr1[219,1] = '16:40'
$arr1[219,2] = '225_GOMEZ'
$arr1[219,3] = 44
$arr1[219,4] = 'LP1912'
$arr1[220,0] = '14:46:12'
$arr1[220,1] = '16:41'
$arr1[220,2] = '225_GOMEZ'
$arr1[220,3] = 115
$arr1[220,4] = 'LP1912'
$arr1[221,0] = '14:53:29'
$arr1[221,1] = '16:46'
$arr1[221,2] = '17_ROMERO'
$arr1[221,3] = 113
$arr1[221,4] = 'LP1912'
$arr1[222,0] = '15:16:46'
$arr1[222,1] = '16:53'
$arr1[222,2] = '11_ETCHEVERRY'
$arr1[222,3] = 97
$arr1[222,4] = 'LP1912'
$arr1[223,0] = '16:12:06'
$arr1[223,1] = '16:54'
$arr1[223,2] = '11_ETCHEVERRY'
$arr1[223,3] = 42
$arr1[223,4] = 'LP1912'
$arr1[224,0] = '16:12:06'
$arr1[224,1] = '16:55'
$arr1[224,2] = '16_SANTA ANA'
$arr1[224,3] = 43
$arr1[224,4] = 'LP1912'
$arr1[225,0] = '15:56:56'
$arr1[225,1] = '16:57'
$arr1[225,2] = '15_ABASTO'
$arr1[225,3] = 61
$arr1[225,4] = 'LP1912'
$arr1[226,0] = '15:16:46'
$arr1[226,1] = '16:58'
$arr1[226,2] = '15_ABASTO'
$arr1[226,3] = 102
$arr1[226,4] = 'LP1912'
$arr1[227,0] = '15:56:56'
$arr1[227,1] = '17:01'
$arr1[227,2] = '23_HERNANDEZ'
$arr1[227,3] = 65
$arr1[227,4] = 'LP1912'
$arr1[228,0] = '15:44:42'
$arr1[228,1] = '17:02'
$arr1[228,2] = '23_HERNANDEZ'
$arr1[228,3] = 78
$arr1[228,4] = 'LP1912'
$arr1[229,0] = '16:28:21'
$arr1[229,1] = '17:06'
$arr1[229,2] = '23_HERNANDEZ'
$arr1[229,3] = 38
$arr1[229,4] = 'LP1912'
$arr1[230,0] = '16:44:58'
$arr1[230,1] = '17:06'
$arr1[230,2] = '16_P MOR-SANTA ANA'
$arr1[230,3] = 22
$arr1[230,4] = 'LP1912'
$arr1[231,0] = '15:16:46'
$arr1[231,1] = '17:07'
$arr1[231,2] = '16_P MOR-SANTA ANA'
$arr1[231,3] = 111
$arr1[231,4] = 'LP1912'
$arr1[232,0] = '16:28:21'
$arr1[232,1] = '17:08'
$arr1[232,2] = '10_OLMOS'
$arr1[232,3] = 40
$arr1[232,4] = 'LP1912'
$arr1[233,0] = '15:16:46'
$arr1[233,1] = '17:09'
$arr1[233,2] = '215C_EL PATO'
$arr1[233,3] = 113
$arr1[233,4] = 'LP1912'
$arr1[234,0] = '16:12:06'
$arr1[234,1] = '17:10'
$arr1[234,2] = '215C_EL PATO'
$arr1[234,3] = 58
$arr1[234,4] = 'LP1912'
$arr1[235,0] = '16:44:58'
$arr1[235,1] = '17:20'
$arr1[235,2] = '15X38_ABASTO'
$arr1[235,3] = 36
$arr1[235,4] = 'LP1912'
$arr1[236,0] = '15:44:42'
$arr1[236,1] = '17:21'
$arr1[236,2] = '15X38_ABASTO'
$arr1[236,3] = 97
$arr1[236,4] = 'LP1912'
$arr1[237,0] = '16:28:21'
$arr1[237,1] = '17:32'
$arr1[237,2] = '27_EL RETIRO'
$arr1[237,3] = 64
$arr1[237,4] = 'LP1912'
$arr1[238,0] = '15:56:56'
$arr1[238,1] = '17:33'
$arr1[238,2] = '17_ROMERO'
$arr1[238,3] = 97
$arr1[238,4] = 'LP1912'
$arr1[239,0] = '15:44:42'
$arr1[239,1] = '17:34'
$arr1[239,2] = '17_ROMERO'
$arr1[239,3] = 110
$arr1[239,4] = 'LP1912'
$arr1[240,0] = '15:44:42'
$arr1[240,1] = '17:36'
$arr1[240,2] = '27_EL RETIRO'
$arr1[240,3] = 112
$arr1[240,4] = 'LP1912'
$arr1[241,0] = '17:35:41'
$arr1[241,1] = '17:37'
$arr1[241,2] = '215B_EL PATO'
$arr1[241,3] = 2
$arr1[241,4] = 'LP1912'
$arr1[242,0] = '16:12:06'
$arr1[242,1] = '17:37'
$arr1[242,2] = '27_EL RETIRO'
$arr1[242,3] = 85
$arr1[242,4] = 'LP1912'
$arr1[243,0] = '15:44:42'
$arr1[243,1] = '17:38'
$arr1[243,2] = '215B_EL PATO'
$arr1[243,3] = 114
$arr1[243,4] = 'LP1912'
$arr1[244,0] = '16:12:06'
$arr1[244,1] = '17:39'
$arr1[244,2] = '215B_EL PATO'
$arr1[244,3] = 87
$arr1[244,4] = 'LP1912'
$arr1[245,0] = '15:56:56'
$arr1[245,1] = '17:45'
$arr1[245,2] = '215_EL PELIGRO'
$arr1[245,3] = 109
$arr1[245,4] = 'LP1912'
$arr1[246,0] = '16:12:06'
$arr1[246,1] = '17:46'
$arr1[246,2] = '215_EL PELIGRO'
$arr1[246,3] = 94
$arr1[246,4] = 'LP1912'
$arr1[247,0] = '17:47:45'
$arr1[247,1] = '17:47'
$arr1[247,2] = '215_EL PELIGRO'
$arr1[247,3] = 0
$arr1[247,4] = 'LP1912'
$arr1[248,0] = '16:12:06'
$arr1[248,1] = '17:49'
$arr1[248,2] = '10_OLMOS'
$arr1[248,3] = 97
$arr1[248,4] = 'LP1912'
$arr1[249,0] = '17:13:30'
$arr1[249,1] = '17:51'
$arr1[249,2] = '23_HERNANDEZ'
$arr1[249,3] = 38
$arr1[249,4] = 'LP1912'
$arr1[250,0] = '15:56:56'
$arr1[250,1] = '17:51'
$arr1[250,2] = '10_OLMOS'
$arr1[250,3] = 115
$arr1[250,4] = 'LP1912'
$arr1[251,0] = '16:28:21'
$arr1[251,1] = '17:52'
$arr1[251,2] = '23_HERNANDEZ'
$arr1[251,3] = 84
$arr1[251,4] = 'LP1912'
$arr1[252,0] = '16:51:51'
$arr1[252,1] = '17:53'
$arr1[252,2] = '10_OLMOS'
$arr1[252,3] = 62
$arr1[252,4] = 'LP1912'
$arr1[253,0] = '16:37:37'
$arr1[253,1] = '17:53'
$arr1[253,2] = '23_HERNANDEZ'
$arr1[253,3] = 76
$arr1[253,4] = 'LP1912'
$arr1[254,0] = '17:55:25'
$arr1[254,1] = '17:55'
$arr1[254,2] = '10_OLMOS'
$arr1[254,3] = 0
$arr1[254,4] = 'LP1912'
$arr1[255,0] = '16:44:58'
$arr1[255,1] = '17:57'
$arr1[255,2] = '17_ROMERO'
$arr1[255,3] = 73
$arr1[255,4] = 'LP1912'
$arr1[256,0] = '16:12:06'
$arr1[256,1] = '17:58'
$arr1[256,2] = '17_ROMERO'
$arr1[256,3] = 106
$arr1[256,4] = 'LP1912'
$arr1[257,0] = '16:28:21'
$arr1[257,1] = '18:05'
$arr1[257,2] = '11_ETCHEVERRY'
$arr1[257,3] = 97
$arr1[257,4] = 'LP1912'
$arr1[258,0] = '16:12:06'
$arr1[258,1] = '18:06'
$arr1[258,2] = '11_ETCHEVERRY'
$arr1[258,3] = 114
$arr1[258,4] = 'LP1912'
$arr1[259,0] = '16:44:58'
$arr1[259,1] = '18:09'
$arr1[259,2] = '16_P MOR-SANTA ANA'
$arr1[259,3] = 85
$arr1[259,4] = 'LP1912'
$arr1[260,0] = '16:44:58'
$arr1[260,1] = '18:09'
$arr1[260,2] = '15_ABASTO'
$arr1[260,3] = 85
$arr1[260,4] = 'LP1912'
$arr1[261,0] = '16:12:06'
$arr1[261,1] = '18:10'
$arr1[261,2] = '15_ABASTO'
$arr1[261,3] = 118
$arr1[261,4] = 'LP1912'
$arr1[262,0] = '16:12:06'
$arr1[262,1] = '18:10'
$arr1[262,2] = '16_P MOR-SANTA ANA'
$arr1[262,3] = 118
$arr1[262,4] = 'LP1912'
$arr1[263,0] = '18:11:09'
$arr1[263,1] = '18:11'
$arr1[263,2] = '15_ABASTO'
$arr1[263,3] = 0
$arr1[263,4] = 'LP1912'
$arr1[264,0] = '18:11:09'
$arr1[264,1] = '18:11'
$arr1[264,2] = '16_P MOR-SANTA ANA'
$arr1[264,3] = 0
$arr1[264,4] = 'LP1912'
$arr1[265,0] = '16:44:58'
$arr1[265,1] = '18:16'
$arr1[265,2] = '10_OLMOS'
$arr1[265,3] = 92
$arr1[265,4] = 'LP1912'
$arr1[266,0] = '16:28:21'
$arr1[266,1] = '18:17'
$arr1[266,2] = '10_OLMOS'
$arr1[266,3] = 109
$arr1[266,4] = 'LP1912'
$arr1[267,0] = '16:37:37'
$arr1[267,1] = '18:21'
$arr1[267,2] = '215C_EL PATO'
$arr1[267,3] = 104
$arr1[267,4] = 'LP1912'
$arr1[268,0] = '16:28:21'
$arr1[268,1] = '18:22'
$arr1[268,2] = '215C_EL PATO'
$arr1[268,3] = 114
$arr1[268,4] = 'LP1912'
$arr1[269,0] = '16:28:21'
$arr1[269,1] = '18:25'
$arr1[269,2] = '16_SANTA ANA'
$arr1[269,3] = 117
$arr1[269,4] = 'LP1912'
$arr1[270,0] = '17:13:30'
$arr1[270,1] = '18:29'
$arr1[270,2] = '23_HERNANDEZ'
$arr1[270,3] = 76
$arr1[270,4] = 'LP1912'
$arr1[271,0] = '16:37:37'
$arr1[271,1] = '18:29'
$arr1[271,2] = '14_ABASTO'
$arr1[271,3] = 112
$arr1[271,4] = 'LP1912'
$arr1[272,0] = '17:55:25'
$arr1[272,1] = '18:30'
$arr1[272,2] = '14_ABASTO'
$arr1[272,3] = 35
$arr1[272,4] = 'LP1912'
$arr1[273,0] = '18:30:48'
$arr1[273,1] = '18:31'
$arr1[273,2] = '14_ABASTO'
$arr1[273,3] = 1
$arr1[273,4] = 'LP1912'
$arr1[274,0] = '17:47:45'
$arr1[274,1] = '18:34'
$arr1[274,2] = '23_HERNANDEZ'
$arr1[274,3] = 47
$arr1[274,4] = 'LP1912'
$arr1[275,0] = '16:44:58'
$arr1[275,1] = '18:35'
$arr1[275,2] = '15X38_ABASTO'
$arr1[275,3] = 111
$arr1[275,4] = 'LP1912'
$arr1[276,0] = '16:37:37'
$arr1[276,1] = '18:36'
$arr1[276,2] = '15X38_ABASTO'
$arr1[276,3] = 119
$arr1[276,4] = 'LP1912'
$arr1[277,0] = '17:35:41'
$arr1[277,1] = '18:37'
$arr1[277,2] = '23_HERNANDEZ'
$arr1[277,3] = 62
$arr1[277,4] = 'LP1912'
$arr1[278,0] = '16:44:58'
$arr1[278,1] = '18:40'
$arr1[278,2] = '10_OLMOS'
$arr1[278,3] = 116
$arr1[278,4] = 'LP1912'
$arr1[279,0] = '17:13:30'
$arr1[279,1] = '18:41'
$arr1[279,2] = '10_OLMOS'
$arr1[279,3] = 88
$arr1[279,4] = 'LP1912'
$arr1[280,0] = '17:13:30'
$arr1[280,1] = '18:45'
$arr1[280,2] = '16_SANTA ANA'
$arr1[280,3] = 92
$arr1[280,4] = 'LP1912'
$arr1[281,0] = '17:55:25'
$arr1[281,1] = '18:46'
$arr1[281,2] = '16_SANTA ANA'
$arr1[281,3] = 51
$arr1[281,4] = 'LP1912'
$arr1[282,0] = '18:30:48'
$arr1[282,1] = '18:48'
$arr1[282,2] = '10_OLMOS'
$arr1[282,3] = 18
$arr1[282,4] = 'LP1912'
$arr1[283,0] = '17:13:30'
$arr1[283,1] = '18:52'
$arr1[283,2] = '17_ROMERO'
$arr1[283,3] = 99
$arr1[283,4] = 'LP1912'
$arr1[284,0] = '17:13:30'
$arr1[284,1] = '18:57'
$arr1[284,2] = '16_P MOR-SANTA ANA'
$arr1[284,3] = 104
$arr1[284,4] = 'LP1912'
$arr1[285,0] = '17:13:30'
$arr1[285,1] = '18:59'
$arr1[285,2] = '14_ABASTO'
$arr1[285,3] = 106
$arr1[285,4] = 'LP1912'
$arr1[286,0] = '18:30:48'
$arr1[286,1] = '19:00'
$arr1[286,2] = '14_ABASTO'
$arr1[286,3] = 30
$arr1[286,4] = 'LP1912'
$arr1[287,0] = '17:47:45'
$arr1[287,1] = '19:02'
$arr1[287,2] = '14_ABASTO'
$arr1[287,3] = 75
$arr1[287,4] = 'LP1912'
$arr1[288,0] = '17:35:41'
$arr1[288,1] = '19:03'
$arr1[288,2] = '215_EL PELIGRO'
$arr1[288,3] = 88
$arr1[288,4] = 'LP1912'
$arr1[289,0] = '17:55:25'
$arr1[289,1] = '19:03'
$arr1[289,2] = '14_ABASTO'
$arr1[289,3] = 68
$arr1[289,4] = 'LP1912'
$arr1[290,0] = '17:13:30'
$arr1[290,1] = '19:04'
$arr1[290,2] = '215_EL PELIGRO'
$arr1[290,3] = 111
$arr1[290,4] = 'LP1912'
$arr1[291,0] = '18:30:48'
$arr1[291,1] = '19:10'
$arr1[291,2] = '16_SANTA ANA'
$arr1[291,3] = 40
$arr1[291,4] = 'LP1912'
$arr1[292,0] = '17:55:25'
$arr1[292,1] = '19:11'
$arr1[292,2] = '16_SANTA ANA'
$arr1[292,3] = 76
$arr1[292,4] = 'LP1912'
$arr1[293,0] = '17:55:25'
$arr1[293,1] = '19:14'
$arr1[293,2] = '27_EL RETIRO'
$arr1[293,3] = 79
$arr1[293,4] = 'LP1912'
$arr1[294,0] = '17:47:45'
$arr1[294,1] = '19:15'
$arr1[294,2] = '17_ROMERO'
$arr1[294,3] = 88
$arr1[294,4] = 'LP1912'
$arr1[295,0] = '17:55:25'
$arr1[295,1] = '19:16'
$arr1[295,2] = '17_ROMERO'
$arr1[295,3] = 81
$arr1[295,4] = 'LP1912'
$arr1[296,0] = '17:35:41'
$arr1[296,1] = '19:16'
$arr1[296,2] = '27_EL RETIRO'
$arr1[296,3] = 101
$arr1[296,4] = 'LP1912'
$arr1[297,0] = '17:35:41'
$arr1[297,1] = '19:17'
$arr1[297,2] = '14X44_ABASTO'
$arr1[297,3] = 102
$arr1[297,4] = 'LP1912'
$arr1[298,0] = '17:55:25'
$arr1[298,1] = '19:22'
$arr1[298,2] = '23_HERNANDEZ'
$arr1[298,3] = 87
$arr1[298,4] = 'LP1912'
$arr1[299,0] = '18:30:48'
$arr1[299,1] = '19:23'
$arr1[299,2] = '16_SANTA ANA'
$arr1[299,3] = 53
$arr1[299,4] = 'LP1912'
$arr1[300,0] = '18:30:48'
$arr1[300,1] = '19:25'
$arr1[300,2] = '23_HERNANDEZ'
$arr1[300,3] = 55
$arr1[300,4] = 'LP1912'
$arr1[301,0] = '18:11:09'
$arr1[301,1] = '19:27'
$arr1[301,2] = '23_HERNANDEZ'
$arr1[301,3] = 76
$arr1[301,4] = 'LP1912'
$arr1[302,0] = '17:35:41'
$arr1[302,1] = '19:27'
$arr1[302,2] = '215C_EL PATO'
$arr1[302,3] = 112
$arr1[302,4] = 'LP1912'
$arr1[303,0] = '17:55:25'
$arr1[303,1] = '19:28'
$arr1[303,2] = '215C_EL PATO'
$arr1[303,3] = 93
$arr1[303,4] = 'LP1912'
$arr1[304,0] = '17:47:45'
$arr1[304,1] = '19:35'
$arr1[304,2] = '11_ETCHEVERRY'
$arr1[304,3] = 108
$arr1[304,4] = 'LP1912'
$arr1[305,0] = '17:55:25'
$arr1[305,1] = '19:36'
$arr1[305,2] = '11_ETCHEVERRY'
$arr1[305,3] = 101
$arr1[305,4] = 'LP1912'
$arr1[306,0] = '17:55:25'
$arr1[306,1] = '19:39'
$arr1[306,2] = '15X38_ABASTO'
$arr1[306,3] = 104
$arr1[306,4] = 'LP1912'
$arr1[307,0] = '17:47:45'
$arr1[307,1] = '19:42'
$arr1[307,2] = '15X38_ABASTO'
$arr1[307,3] = 115
$arr1[307,4] = 'LP1912'
$arr1[308,0] = '17:55:25'
$arr1[308,1] = '19:52'
$arr1[308,2] = '81_EL PELIGRO'
$arr1[308,3] = 117
$arr1[308,4] = 'LP1912'
$arr1[309,0] = '17:55:25'
$arr1[309,1] = '19:53'
$arr1[309,2] = '225_GOMEZ'
$arr1[309,3] = 118
$arr1[309,4] = 'LP1912'
$arr1[310,0] = '18:11:09'
$arr1[310,1] = '20:06'
$arr1[310,2] = '215C_EL PATO'
$arr1[310,3] = 115
$arr1[310,4] = 'LP1912'
$arr1[311,0] = '18:30:48'
$arr1[311,1] = '20:21'
$arr1[311,2] = '15_ABASTO'
$arr1[311,3] = 111
$arr1[311,4] = 'LP1912'
$ws1.Range("A6:E317").Value = $arr1

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = 'Última actualización: 18:30:48'

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = 'Última actualización: 18:30:48'
$ws3.Range("A3").Value = 'Total filas: 45'
$arr3 = New-Object 'object[,]' 45,5
$arr3[0,0] = '06:16:41'
$arr3[0,1] = '07:42'
$arr3[0,2] = '215A_LA PLATA'
$arr3[0,3] = 86
$arr3[0,4] = 'L6173'
$arr3[1,0] = '05:57:04'
$arr3[1,1] = '07:43'
$arr3[1,2] = '215A_LA PLATA'
$arr3[1,3] = 106
$arr3[1,4] = 'L6173'
$arr3[2,0] = '06:53:44'
$arr3[2,1] = '08:35'
$arr3[2,2] = '215A_LA PLATA'
$arr3[2,3] = 102
$arr3[2,4] = 'L6173'
$arr3[3,0] = '06:46:20'
$arr3[3,1] = '08:36'
$arr3[3,2] = '215A_LA PLATA'
$arr3[3,3] = 110
$arr3[3,4] = 'L6173'
$arr3[4,0] = '08:38:24'
$arr3[4,1] = '08:38'
$arr3[4,2] = '215A_LA PLATA'
$arr3[4,3] = 0
$arr3[4,4] = 'L6173'
$arr3[5,0] = '06:53:44'
$arr3[5,1] = '08:50'
$arr3[5,2] = '215C_LA PLATA'
$arr3[5,3] = 117
$arr3[5,4] = 'L6203'
$arr3[6,0] = '07:13:03'
$arr3[6,1] = '08:51'
$arr3[6,2] = '215C_LA PLATA'
$arr3[6,3] = 98
$arr3[6,4] = 'L6203'
$arr3[7,0] = '08:52:40'
$arr3[7,1] = '08:55'
$arr3[7,2] = '215C_LA PLATA'
$arr3[7,3] = 3
$arr3[7,4] = 'L6203'
$arr3[8,0] = '08:28:52'
$arr3[8,1] = '09:20'
$arr3[8,2] = '215A_LA PLATA'
$arr3[8,3] = 52
$arr3[8,4] = 'L6173'
$arr3[9,0] = '07:49:32'
$arr3[9,1] = '09:21'
$arr3[9,2] = '215A_LA PLATA'
$arr3[9,3] = 92
$arr3[9,4] = 'L6173'
$arr3[10,0] = '09:22:34'
$arr3[10,1] = '09:23'
$arr3[10,2] = '215A_LA PLATA'
$arr3[10,3] = 1
$arr3[10,4] = 'L6173'
$arr3[11,0] = '08:28:52'
$arr3[11,1] = '10:12'
$arr3[11,2] = '215C_LA PLATA'
$arr3[11,3] = 104
$arr3[11,4] = 'L6203'
$arr3[12,0] = '08:38:24'
$arr3[12,1] = '10:13'
$arr3[12,2] = '215C_LA PLATA'
$arr3[12,3] = 95
$arr3[12,4] = 'L6203'
$arr3[13,0] = '08:52:40'
$arr3[13,1] = '10:29'
$arr3[13,2] = '215B_LP-P MOR-1 Y 57'
$arr3[13,3] = 97
$arr3[13,4] = 'L6173'
$arr3[14,0] = '08:38:24'
$arr3[14,1] = '10:30'
$arr3[14,2] = '215B_LP-P MOR-1 Y 57'
$arr3[14,3] = 112
$arr3[14,4] = 'L6173'
$arr3[15,0] = '08:52:40'
$arr3[15,1] = '10:30'
$arr3[15,2] = '215A_LA PLATA'
$arr3[15,3] = 98
$arr3[15,4] = 'L6173'
$arr3[16,0] = '08:45:31'
$arr3[16,1] = '10:31'
$arr3[16,2] = '215A_LA PLATA'
$arr3[16,3] = 106
$arr3[16,4] = 'L6173'
$arr3[17,0] = '10:36:50'
$arr3[17,1] = '11:25'
$arr3[17,2] = '215C_LA PLATA'
$arr3[17,3] = 49
$arr3[17,4] = 'L6203'
$arr3[18,0] = '10:04:30'
$arr3[18,1] = '11:26'
$arr3[18,2] = '215C_LA PLATA'
$arr3[18,3] = 82
$arr3[18,4] = 'L6203'
$arr3[19,0] = '11:33:52'
$arr3[19,1] = '13:11'
$arr3[19,2] = '215C_LA PLATA'
$arr3[19,3] = 98
$arr3[19,4] = 'L6203'
$arr3[20,0] = '11:13:15'
$arr3[20,1] = '13:12'
$arr3[20,2] = '215C_LA PLATA'
$arr3[20,3] = 119
$arr3[20,4] = 'L6203'
$arr3[21,0] = '13:14:31'
$arr3[21,1] = '13:16'
$arr3[21,2] = '215C_LA PLATA'
$arr3[21,3] = 2
$arr3[21,4] = 'L6203'
$arr3[22,0] = '11:33:52'
$arr3[22,1] = '13:20'
$arr3[22,2] = '215B_LP-P MOR-1 Y 57'
$arr3[22,3] = 107
$arr3[22,4] = 'L6173'
$arr3[23,0] = '11:46:32'
$arr3[23,1] = '13:21'
$arr3[23,2] = '215B_LP-P MOR-1 Y 57'
$arr3[23,3] = 95
$arr3[23,4] = 'L6173'
$arr3[24,0] = '12:11:21'
$arr3[24,1] = '13:57'
$arr3[24,2] = '215C_LA PLATA'
$arr3[24,3] = 106
$arr3[24,4] = 'L6203'
$arr3[25,0] = '13:55:43'
$arr3[25,1] = '13:58'
$arr3[25,2] = '215C_LA PLATA'
$arr3[25,3] = 3
$arr3[25,4] = 'L6203'
$arr3[26,0] = '13:14:31'
$arr3[26,1] = '14:03'
$arr3[26,2] = '215C_LA PLATA'
$arr3[26,3] = 49
$arr3[26,4] = 'L6203'
$arr3[27,0] = '13:55:43'
$arr3[27,1] = '14:26'
$arr3[27,2] = '215C_LA PLATA'
$arr3[27,3] = 31
$arr3[27,4] = 'L6203'
$arr3[28,0] = '12:46:07'
$arr3[28,1] = '14:27'
$arr3[28,2] = '215C_LA PLATA'
$arr3[28,3] = 101
$arr3[28,4] = 'L6203'
$arr3[29,0] = '15:16:46'
$arr3[29,1] = '15:19'
$arr3[29,2] = '215A_LA PLATA'
$arr3[29,3] = 3
$arr3[29,4] = 'L6173'
$arr3[30,0] = '13:55:43'
$arr3[30,1] = '15:21'
$arr3[30,2] = '215A_LA PLATA'
$arr3[30,3] = 86
$arr3[30,4] = 'L6173'
$arr3[31,0] = '13:41:21'
$arr3[31,1] = '15:22'
$arr3[31,2] = '215A_LA PLATA'
$arr3[31,3] = 101
$arr3[31,4] = 'L6173'
$arr3[32,0] = '14:32:44'
$arr3[32,1] = '16:01'
$arr3[32,2] = '215C_LA PLATA'
$arr3[32,3] = 89
$arr3[32,4] = 'L6203'
$arr3[33,0] = '14:11:28'
$arr3[33,1] = '16:02'
$arr3[33,2] = '215C_LA PLATA'
$arr3[33,3] = 111
$arr3[33,4] = 'L6203'
$arr3[34,0] = '14:32:44'
$arr3[34,1] = '16:29'
$arr3[34,2] = '215B_LP-P MOR-40 Y 115'
$arr3[34,3] = 117
$arr3[34,4] = 'L6173'
$arr3[35,0] = '14:46:12'
$arr3[35,1] = '16:30'
$arr3[35,2] = '215B_LP-P MOR-40 Y 115'
$arr3[35,3] = 104
$arr3[35,4] = 'L6173'
$arr3[36,0] = '16:28:21'
$arr3[36,1] = '16:31'
$arr3[36,2] = '215B_LP-P MOR-40 Y 115'
$arr3[36,3] = 3
$arr3[36,4] = 'L6173'
$arr3[37,0] = '15:16:46'
$arr3[37,1] = '17:05'
$arr3[37,2] = '215C_LA PLATA'
$arr3[37,3] = 109
$arr3[37,4] = 'L6203'
$arr3[38,0] = '16:12:06'
$arr3[38,1] = '17:06'
$arr3[38,2] = '215C_LA PLATA'
$arr3[38,3] = 54
$arr3[38,4] = 'L6203'
$arr3[39,0] = '16:37:37'
$arr3[39,1] = '18:35'
$arr3[39,2] = '215C_LA PLATA'
$arr3[39,3] = 118
$arr3[39,4] = 'L6203'
$arr3[40,0] = '17:13:30'
$arr3[40,1] = '18:36'
$arr3[40,2] = '215C_LA PLATA'
$arr3[40,3] = 83
$arr3[40,4] = 'L6203'
$arr3[41,0] = '17:35:41'
$arr3[41,1] = '19:23'
$arr3[41,2] = '215B_LP-P MOR-1 Y 57'
$arr3[41,3] = 108
$arr3[41,4] = 'L6173'
$arr3[42,0] = '17:55:25'
$arr3[42,1] = '19:24'
$arr3[42,2] = '215B_LP-P MOR-1 Y 57'
$arr3[42,3] = 89
$arr3[42,4] = 'L6173'
$arr3[43,0] = '18:30:48'
$arr3[43,1] = '19:57'
$arr3[43,2] = '215C_LA PLATA'
$arr3[43,3] = 87
$arr3[43,4] = 'L6203'
$arr3[44,0] = '18:11:09'
$arr3[44,1] = '19:58'
$arr3[44,2] = '215C_LA PLATA'
$arr3[44,3] = 107
$arr3[44,4] = 'L6203'
$ws3.Range("A6:E50").Value = $arr3
